$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(124, 8).Value = 20820
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 20820
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 20820
$ws.Cells.Item(124, 14).Value = -30640
$ws.Cells.Item(125, 8).Value = 6979.385
$ws.Cells.Item(125, 9).Value = 6433
$ws.Cells.Item(125, 10).Value = 7222.222
$ws.Cells.Item(125, 11).Value = 57897
$ws.Cells.Item(125, 12).Value = 64999.998
$ws.Cells.Item(125, 13).Value = -55437
$ws.Cells.Item(125, 14).Value = -69919.99799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 716.7143
$ws.Cells.Item(110, 9).Value = 621.0476
$ws.Cells.Item(110, 10).Value = 1003.7143
$ws.Cells.Item(110, 11).Value = 621.0476
$ws.Cells.Item(110, 12).Value = 1003.7143
$ws.Cells.Item(110, 13).Value = 1423.9524
$ws.Cells.Item(110, 14).Value = -5093.7143
$ws.Cells.Item(122, 8).Value = 1071255
$ws.Cells.Item(122, 9).Value = 1223724.9
$ws.Cells.Item(122, 10).Value = 3966.3333
$ws.Cells.Item(122, 11).Value = 3671174.7
$ws.Cells.Item(122, 12).Value = 11898.9999
$ws.Cells.Item(122, 13).Value = -3668724.7
$ws.Cells.Item(122, 14).Value = -16798.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 16611.705
$ws.Cells.Item(20, 9).Value = 1766.5555
$ws.Cells.Item(20, 10).Value = 33312.5
$ws.Cells.Item(20, 11).Value = 1766.5555
$ws.Cells.Item(20, 12).Value = 33312.5
$ws.Cells.Item(20, 13).Value = -1519.5555
$ws.Cells.Item(20, 14).Value = -33806.5
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).Value = $null
$ws.Cells.Item(107, 8).Value = 1096.1765
$ws.Cells.Item(107, 9).Value = 1118.5
$ws.Cells.Item(107, 10).Value = 1042.6
$ws.Cells.Item(107, 11).Value = 1118.5
$ws.Cells.Item(107, 12).Value = 1042.6
$ws.Cells.Item(107, 13).Value = 801.5
$ws.Cells.Item(107, 14).Value = -4882.6
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).Value = $null
$ws.Cells.Item(133, 8).Value = 30890
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 30890
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 30890
$ws.Cells.Item(133, 14).Value = -41010
$ws.Cells.Item(134, 8).Value = 5864.074
$ws.Cells.Item(134, 9).Value = 6515
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 19545
$ws.Cells.Item(134, 12).Value = 9000
$ws.Cells.Item(134, 13).Value = -17010
$ws.Cells.Item(134, 14).Value = -14070
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).Value = $null
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).Value = $null
$ws.Cells.Item(137, 8).Value = 39709
$ws.Cells.Item(137, 9).Value = 39709
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 39709
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -34609
$ws.Cells.Item(138, 8).Value = 46240
$ws.Cells.Item(138, 9).Value = 2000
$ws.Cells.Item(138, 10).Value = 57300
$ws.Cells.Item(138, 11).Value = 2000
$ws.Cells.Item(138, 12).Value = 57300
$ws.Cells.Item(138, 13).Value = 3140
$ws.Cells.Item(138, 14).Value = -67580
$ws.Cells.Item(139, 8).Value = 65854.5
$ws.Cells.Item(139, 9).Value = 20709
$ws.Cells.Item(139, 10).Value = 111000
$ws.Cells.Item(139, 11).Value = 20709
$ws.Cells.Item(139, 12).Value = 111000
$ws.Cells.Item(139, 13).Value = -15569
$ws.Cells.Item(139, 14).Value = -121280
$ws.Cells.Item(140, 8).Value = 49237.375
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 49237.375
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 49237.375
$ws.Cells.Item(140, 14).Value = -59597.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 7693786.5
$ws.Cells.Item(16, 9).Value = 10990352
$ws.Cells.Item(16, 10).Value = 1800
$ws.Cells.Item(16, 11).Value = 10990352
$ws.Cells.Item(16, 12).Value = 1800
$ws.Cells.Item(16, 13).Value = -10990065
$ws.Cells.Item(16, 14).Value = -2374
$ws.Cells.Item(31, 8).Value = 7826.357
$ws.Cells.Item(31, 9).Value = 2579.9167
$ws.Cells.Item(31, 10).Value = 11761.1875
$ws.Cells.Item(31, 11).Value = 2579.9167
$ws.Cells.Item(31, 12).Value = 11761.1875
$ws.Cells.Item(31, 13).Value = -2284.9167
$ws.Cells.Item(31, 14).Value = -12351.1875
$ws.Cells.Item(34, 8).Value = 7826.357
$ws.Cells.Item(34, 9).Value = 2579.9167
$ws.Cells.Item(34, 10).Value = 11761.1875
$ws.Cells.Item(34, 11).Value = 2579.9167
$ws.Cells.Item(34, 12).Value = 11761.1875
$ws.Cells.Item(34, 13).Value = -2377.9167
$ws.Cells.Item(34, 14).Value = -12165.1875
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = $null
$ws.Cells.Item(64, 14).Value = $null
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = $null
$ws.Cells.Item(67, 14).Value = $null
$ws.Cells.Item(113, 8).Value = 7693786.5
$ws.Cells.Item(113, 9).Value = 10990352
$ws.Cells.Item(113, 10).Value = 1800
$ws.Cells.Item(113, 11).Value = 10990352
$ws.Cells.Item(113, 12).Value = 1800
$ws.Cells.Item(113, 13).Value = -10988182
$ws.Cells.Item(113, 14).Value = -6140
$ws.Cells.Item(134, 8).Value = 3275.12
$ws.Cells.Item(134, 9).Value = 3707.158
$ws.Cells.Item(134, 10).Value = 1907
$ws.Cells.Item(134, 11).Value = 11121.474
$ws.Cells.Item(134, 12).Value = 5721
$ws.Cells.Item(134, 13).Value = -8586.474
$ws.Cells.Item(134, 14).Value = -10791
$ws.Cells.Item(138, 8).Value = 63200
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 63200
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 63200
$ws.Cells.Item(138, 14).Value = -73480
$ws.Cells.Item(140, 8).Value = 28454.85
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 28454.85
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 28454.85
$ws.Cells.Item(140, 14).Value = -38814.85
$ws.Cells.Item(141, 8).Value = 21066.334
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 21066.334
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 21066.334
$ws.Cells.Item(141, 14).Value = -31426.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 2148.4
$ws.Cells.Item(103, 9).Value = 394.66666
$ws.Cells.Item(103, 10).Value = 2900
$ws.Cells.Item(103, 11).Value = 1183.99998
$ws.Cells.Item(103, 12).Value = 8700
$ws.Cells.Item(103, 13).Value = -304.9999800000001
$ws.Cells.Item(103, 14).Value = -10458
$ws.Cells.Item(130, 8).Value = 6583.3335
$ws.Cells.Item(130, 9).Value = 1000
$ws.Cells.Item(130, 10).Value = 7700
$ws.Cells.Item(130, 11).Value = 3000
$ws.Cells.Item(130, 12).Value = 23100
$ws.Cells.Item(130, 13).Value = 2020
$ws.Cells.Item(130, 14).Value = -33140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 35055
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 35055
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 35055
$ws.Cells.Item(63, 14).Value = -36427
$ws.Cells.Item(66, 8).Value = 35055
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 35055
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 105165
$ws.Cells.Item(66, 14).Value = -112029
$ws.Cells.Item(68, 8).Value = 42300
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 42300
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 42300
$ws.Cells.Item(68, 14).Value = -43922
$ws.Cells.Item(71, 8).Value = 42300
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 42300
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 126900
$ws.Cells.Item(71, 14).Value = -135012
$ws.Cells.Item(74, 8).Value = 39946.95
$ws.Cells.Item(74, 9).Value = 39999
$ws.Cells.Item(74, 10).Value = 39909.09
$ws.Cells.Item(74, 11).Value = 39999
$ws.Cells.Item(74, 12).Value = 39909.09
$ws.Cells.Item(74, 13).Value = -39063
$ws.Cells.Item(74, 14).Value = -41781.09
$ws.Cells.Item(77, 8).Value = 39946.95
$ws.Cells.Item(77, 9).Value = 39999
$ws.Cells.Item(77, 10).Value = 39909.09
$ws.Cells.Item(77, 11).Value = 119997
$ws.Cells.Item(77, 12).Value = 119727.27
$ws.Cells.Item(77, 13).Value = -115317
$ws.Cells.Item(77, 14).Value = -129087.27
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).Value = $null
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).Value = $null
$ws.Cells.Item(86, 8).Value = 30000
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 30000
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 30000
$ws.Cells.Item(86, 14).Value = -32372
$ws.Cells.Item(89, 8).Value = 30000
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 30000
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 90000
$ws.Cells.Item(89, 14).Value = -101856
$ws.Cells.Item(113, 8).Value = 90910480
$ws.Cells.Item(113, 9).Value = 166667660
$ws.Cells.Item(113, 10).Value = 1860
$ws.Cells.Item(113, 11).Value = 166667660
$ws.Cells.Item(113, 12).Value = 1860
$ws.Cells.Item(113, 13).Value = -166665490
$ws.Cells.Item(113, 14).Value = -6200
$ws.Cells.Item(122, 8).Value = 8833893
$ws.Cells.Item(122, 9).Value = 10804268
$ws.Cells.Item(122, 10).Value = 7145000
$ws.Cells.Item(122, 11).Value = 32412804
$ws.Cells.Item(122, 12).Value = 21435000
$ws.Cells.Item(122, 13).Value = -32410354
$ws.Cells.Item(122, 14).Value = -21439900

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3703606.8
$ws.Cells.Item(122, 9).Value = 4204109
$ws.Cells.Item(122, 10).Value = 2001900
$ws.Cells.Item(122, 11).Value = 12612327
$ws.Cells.Item(122, 12).Value = 6005700
$ws.Cells.Item(122, 13).Value = -12609877
$ws.Cells.Item(122, 14).Value = -6010600
$ws.Cells.Item(132, 8).Value = 18062870
$ws.Cells.Item(132, 9).Value = 24082118
$ws.Cells.Item(132, 10).Value = 5123.8335
$ws.Cells.Item(132, 11).Value = 72246354
$ws.Cells.Item(132, 12).Value = 15371.5005
$ws.Cells.Item(132, 13).Value = -72243824
$ws.Cells.Item(132, 14).Value = -20431.5005
$ws.Cells.Item(136, 8).Value = 10317.12
$ws.Cells.Item(136, 9).Value = 12819.8
$ws.Cells.Item(136, 10).Value = 8648.666999999999
$ws.Cells.Item(136, 11).Value = 38459.39999999999
$ws.Cells.Item(136, 12).Value = 25946.001
$ws.Cells.Item(136, 13).Value = -35909.39999999999
$ws.Cells.Item(136, 14).Value = -31046.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 14).Value = $null
$ws.Cells.Item(100, 8).Value = 8095.5386
$ws.Cells.Item(100, 9).Value = 9449.546
$ws.Cells.Item(100, 10).Value = 648.5
$ws.Cells.Item(100, 11).Value = 18899.092
$ws.Cells.Item(100, 12).Value = 1297
$ws.Cells.Item(100, 13).Value = -18358.092
$ws.Cells.Item(100, 14).Value = -2379
$ws.Cells.Item(122, 8).Value = 1699.25
$ws.Cells.Item(122, 9).Value = 1640.6666
$ws.Cells.Item(122, 10).Value = 1875
$ws.Cells.Item(122, 11).Value = 4921.9998
$ws.Cells.Item(122, 12).Value = 5625
$ws.Cells.Item(122, 13).Value = -2471.9998
$ws.Cells.Item(122, 14).Value = -10525
$ws.Cells.Item(139, 8).Value = 52422.8
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 52422.8
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 52422.8
$ws.Cells.Item(139, 14).Value = -62702.8
